$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FE")

# Expand the Excel Table (ListObject) by 4 rows so it covers A1:M10,
# matching the new data being appended below the existing 6 rows.
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null

# --- Controlled write order so new shared strings land at the same
# --- indices as the target workbook (81..86, in this exact sequence):
#   81 "Todo TRUE - max y min"
#   82 "624gb ram insuficientes"
#   83 "FE05"
#   84 "FE06"
#   85 "FE07"
#   86 "FE08"

# 1) Introduces "Todo TRUE - max y min" (-> shared string 81)
$ws.Cells.Item(7, 11).Value = "Todo TRUE - max y min"   # K7

# 2) Introduces "624gb ram insuficientes" (-> shared string 82), used on
#    the observation column (M) of the four pre-existing data rows.
$ws.Cells.Item(3, 13).Value = "624gb ram insuficientes" # M3
$ws.Cells.Item(4, 13).Value = "624gb ram insuficientes" # M4
$ws.Cells.Item(5, 13).Value = "624gb ram insuficientes" # M5
$ws.Cells.Item(6, 13).Value = "624gb ram insuficientes" # M6

# 3) Introduces "FE05" (-> 83), "FE06" (-> 84), "FE07" (-> 85), "FE08" (-> 86)
$ws.Cells.Item(7, 1).Value  = "FE05"   # A7
$ws.Cells.Item(8, 1).Value  = "FE06"   # A8
$ws.Cells.Item(9, 1).Value  = "FE07"   # A9
$ws.Cells.Item(10, 1).Value = "FE08"   # A10

# --- Fill in the remaining cells for the four new rows ---

# Row 7: FE05
$ws.Cells.Item(7, 2).Value  = 100      # B7 num.trees
$ws.Cells.Item(7, 3).Value  = 5        # C7 max.depth
$ws.Cells.Item(7, 4).Value  = 600      # D7 min.node.size
$ws.Cells.Item(7, 5).Value  = 50       # E7 mtry
$ws.Cells.Item(7, 6).Value  = 10881    # F7 semilla
$ws.Cells.Item(7, 7).Value  = 0.999    # G7 Canarios
$ws.Cells.Item(7, 8).Value  = $true    # H7 lag1
$ws.Cells.Item(7, 9).Value  = $true    # I7 lag2
$ws.Cells.Item(7, 10).Value = $true    # J7 lag3
$ws.Cells.Item(7, 12).Value = "Todo TRUE - max y min"   # L7 Tendencias2

# Row 8: FE06
$ws.Cells.Item(8, 2).Value  = 100      # B8
$ws.Cells.Item(8, 3).Value  = 4        # C8
$ws.Cells.Item(8, 4).Value  = 600      # D8
$ws.Cells.Item(8, 5).Value  = 50       # E8
$ws.Cells.Item(8, 6).Value  = 10881    # F8
$ws.Cells.Item(8, 7).Value  = 0.999    # G8
$ws.Cells.Item(8, 8).Value  = $true    # H8
$ws.Cells.Item(8, 9).Value  = $true    # I8
$ws.Cells.Item(8, 10).Value = $true    # J8
$ws.Cells.Item(8, 11).Value = "Todo TRUE"   # K8
$ws.Cells.Item(8, 12).Value = "Todo TRUE"   # L8

# Row 9: FE07
$ws.Cells.Item(9, 2).Value  = 80       # B9
$ws.Cells.Item(9, 3).Value  = 4        # C9
$ws.Cells.Item(9, 4).Value  = 600      # D9
$ws.Cells.Item(9, 5).Value  = 50       # E9
$ws.Cells.Item(9, 6).Value  = 10881    # F9
$ws.Cells.Item(9, 7).Value  = 0.999    # G9
$ws.Cells.Item(9, 8).Value  = $true    # H9
$ws.Cells.Item(9, 9).Value  = $true    # I9
$ws.Cells.Item(9, 10).Value = $true    # J9
$ws.Cells.Item(9, 11).Value = "Todo TRUE - max y min"   # K9
$ws.Cells.Item(9, 12).Value = "Todo TRUE - max y min"   # L9

# Row 10: FE08
$ws.Cells.Item(10, 2).Value  = 80      # B10
$ws.Cells.Item(10, 3).Value  = 8       # C10
$ws.Cells.Item(10, 4).Value  = 500     # D10
$ws.Cells.Item(10, 5).Value  = 40      # E10
$ws.Cells.Item(10, 6).Value  = 10881   # F10
$ws.Cells.Item(10, 7).Value  = 0.8     # G10
$ws.Cells.Item(10, 8).Value  = $true   # H10
$ws.Cells.Item(10, 9).Value  = $true   # I10
$ws.Cells.Item(10, 10).Value = $true   # J10
$ws.Cells.Item(10, 11).Value = "Todo TRUE - max y min"   # K10
$ws.Cells.Item(10, 12).Value = $false  # L10

# Move the active selection to match the saved cursor position.
$ws.Range("L20").Select() | Out-Null
